$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new column before column B, shifting all subsequent columns right
$ws.Columns("B:B").Insert()

# Set header and value for the new column
$ws.Range("B1").Value = "manufacturerNumber"
$ws.Range("B2").Value = "Hersteller-Artikelnummer"

# Reset selection / view
$ws.Range("A2").Select()
